$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New weight-log rows: date serial, weight (kg), distance (km)
$data = @(
    @(42072, 85.3, 0),
    @(42073, 85.6, 5.31),
    @(42074, 85.6, 7.11),
    @(42075, 85.4, 0),
    @(42076, 85.3, 12.04),
    @(42077, 84.6, 0),
    @(42078, 85.2, 30.13),
    @(42079, 85.2, 0),
    @(42080, 86.4, 0),
    @(42081, 85.2, 0),
    @(42082, 85.7, 0),
    @(42083, 86.1, 0),
    @(42084, 84.6, 12.01),
    @(42085, 85.8, 0),
    @(42086, 84.3, 6.4),
    @(42087, 85, 0),
    @(42088, 84.6, 12.14),
    @(42089, 84.6, 0),
    @(42090, 86, 0),
    @(42091, 85.7, 0),
    @(42092, 85.4, 0),
    @(42093, 85.4, 0),
    @(42094, 85.6, 0),
    @(42095, 85.4, 0),
    @(42096, 84.7, 5.31),
    @(42097, 84.3, 0),
    @(42098, 83.9, 32.14),
    @(42099, 83, 0),
    @(42100, 85.3, 0),
    @(42101, 85.3, 7.52),
    @(42102, 84.5, 0),
    @(42103, 83.5, 0),
    @(42104, 84.3, 0),
    @(42105, 85.2, 0),
    @(42106, 84.5, 12.02),
    @(42107, 85.1, 5),
    @(42108, 84.4, 0),
    @(42109, 84.8, 0),
    @(42110, 84, 11.34),
    @(42111, 83.2, 0),
    @(42112, 84, 0),
    @(42113, 83, 18.04),
    @(42114, 84.1, 10.199999999999999),
    @(42115, 83.9, 0)
)

$firstNewRow = 140
$lastNewRow = $firstNewRow + $data.Count - 1

# Copy the formatting (number format / alignment) of the last existing
# data row (139) down across the whole new block in one shot, then fill
# in the actual values.
$ws.Range("A139:C139").Copy() | Out-Null
$ws.Range("A" + $firstNewRow + ":C" + $lastNewRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$row = $firstNewRow
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}

# Extend the line chart's single series to cover the new rows.
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Sheet1!`$B`$1,Sheet1!`$A`$2:`$A`$" + $lastNewRow + ",Sheet1!`$B`$2:`$B`$" + $lastNewRow + ",1)"

# Match the author's final selection / scroll position.
$ws.Range("C" + $lastNewRow).Select() | Out-Null
try { $excel.ActiveWindow.TopLeftCell = $ws.Range("D28") } catch { }
try { $wb.Windows.Item(1).ScrollRow = 28; $wb.Windows.Item(1).ScrollColumn = 4 } catch { }
